$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "This is good, very good, I liked it, very nice. Really appreciate."

$ws.Range("B3").Select()
